# Updated cryptos list with GitHub Actions: refresh Price / Volume(1h)
# columns for each coin row, plus two pairs of rows (17/18 and 48/49) that
# swapped rank order (their Coin/Link/Price/Volume values moved together).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.094.28'
$ws.Range('E2').Value = '  +6.82%  '
$ws.Range('D3').Value = '3.568.82'
$ws.Range('E3').Value = '  +10.91%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '188.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '552.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.18%  '
$ws.Range('D7').Value = '3.559.77'
$ws.Range('E7').Value = '  +10.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.609'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.70%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  +5.37%  '
$ws.Range('E11').Value = '  +15.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.87'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.73%  '
$ws.Range('E13').Value = '  +8.10%  '
$ws.Range('E14').Value = '  +4.01%  '
$ws.Range('D15').Value = '4.135.41'
$ws.Range('E15').Value = '  +10.89%  '
$ws.Range('D16').Value = '3.573.62'
$ws.Range('E16').Value = '  +10.96%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.121'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.24%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '67.143.41'
$ws.Range('E18').Value = '  +7.13%  '
$ws.Range('E19').Value = '  +6.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.994'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '432.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +18.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '85.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.71%  '
$ws.Range('E25').Value = '  +5.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('E27').Value = '  +10.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.87%  '
$ws.Range('E30').Value = '  +11.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '649.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.18%  '
$ws.Range('E34').Value = '  +4.81%  '
$ws.Range('E35').Value = '  +6.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.59'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.152'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +24.51%  '
$ws.Range('E38').Value = '  +17.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.65'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('E41').Value = '  +4.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.36'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +15.49%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '3.046.46'
$ws.Range('E44').Value = '  +6.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.36%  '
$ws.Range('E47').Value = '  +12.18%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0420'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.90%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.85%  '
$ws.Range('E50').Value = '  +5.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +14.10%  '
